$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 4-17 to reflect repulled data / mean calc
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -4
